# Weekly update: a new Papa (potato) price observation was added to the
# dataset. In the source table, rows are kept in reverse-chronological
# (most-recent-first) order, so the new observation is inserted right above
# the current row 199, shifting all following rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 199 (pushes old rows 199..250 down to 200..251)
$ws.Rows.Item(199).Insert()

# Populate the new row 199 with the new observation
$ws.Cells.Item(199, 1).Value = 7
$ws.Cells.Item(199, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(199, 3).Value = "Ñuble"
$ws.Cells.Item(199, 4).Value = 44551
$ws.Cells.Item(199, 5).Value = 16
$ws.Cells.Item(199, 6).Value = 100114001
$ws.Cells.Item(199, 7).Value = "Papa"
$ws.Cells.Item(199, 8).Value = "Asterix"
$ws.Cells.Item(199, 9).Value = "1a nueva(o)"
$ws.Cells.Item(199, 10).Value = 1000
$ws.Cells.Item(199, 11).Value = 9500
$ws.Cells.Item(199, 12).Value = 10000
$ws.Cells.Item(199, 13).Value = 9750
$ws.Cells.Item(199, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(199, 15).Value = "Región del Maule"
$ws.Cells.Item(199, 16).Value = 390
$ws.Cells.Item(199, 17).Value = 25
$ws.Cells.Item(199, 18).Value = "Hortaliza"
